# Regenerate the s_vals sheet data (filter save games) for holderman_colin.
# Updates columns B-E (and derived sum in G) for rows 2-5.
# Column F (Win flag) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128; D = 0.1575252929769615; E = 0.496779210170732;  G = 5.553084769722144 }
    3 = @{ B = 0.3048080303191223; C = 1.667794583268128; D = 26.21740644021617;  E = 0.496779210170732;  G = 28.68678826397415 }
    4 = @{ B = 0.6753301551942219; C = 0.3127903958511391; D = 26.21740644021617; E = 8.660232485948974;  G = 35.8657594772105 }
    5 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 3.900430680208489; E = 0.496779210170732;  G = 6.740334628841572 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}

$wb.Save()
